$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DA price values for B2:B25
$prices = @(50.81, 51.51, 51.95, 51.91, 52.03, 52.52, 54.46, 58.85, 59.94, 59.95, 59.94, 57.98, 54.11, 56.6, 66.72, 72.21, 78.28, 89.4, 93.73, 87.19, 77.49, 71.62, 70.06, 66.39)

# Capture the number format already used by the (now unused) column E cells
# before we clear them, so it can be reapplied to the B column prices.
$priceFormat = $ws.Range("E2").NumberFormat

for ($i = 0; $i -lt $prices.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $prices[$i]
}

# Apply the numeric style (numFmtId 4 -> "#,##0.00") to B2:B25
$ws.Range("B2:B25").NumberFormat = $priceFormat

# Clear column E (E2:E25 had style-only cells with no value)
$ws.Range("E2:E25").Clear()

# Update selection to B1
$ws.Range("B1").Select()
